$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header dates (row 5) ---
$ws.Range("A5").Value = "Período: 20/10/2025 até 09/11/2025"

# --- Data rows 8-10 (conhecimentos) ---
# Row 8
$ws.Range("A8").Value = "8809, "
$ws.Range("B8").Value = "'3755908"
$ws.Range("C8").Value = " - "
$ws.Range("D8").Value = "'10285063000276"
$ws.Range("E8").Value = "MJ COMERCIO E SERVICOS DE INFORMATICA E TELECOMUNICACOES LTD"
$ws.Range("F8").Value = "'2952192000161"
$ws.Range("G8").Value = "CABO SERVICOS DE TELECOMUNICACOES LTDA"
$ws.Range("H8").Value = "'05/11/2025"
$ws.Range("I8").Value = "DIGITADO"
$ws.Range("J8").Value = "'11/11/2025"
$ws.Range("K8").Value = "'11/11/2025"
$ws.Range("L8").Value = "FINALIZADO"
$ws.Range("M8").Value = "ENTREGA REALIZADA"
$ws.Range("N8").Value = "'08/11/2025"
$ws.Range("O8").Value = "RODO"
$ws.Range("P8").Value = "FOR"
$ws.Range("Q8").Value = "EUSEBIO"
$ws.Range("R8").Value = "CE"
$ws.Range("S8").Value = "NAT"
$ws.Range("T8").Value = "NATAL"
$ws.Range("U8").Value = "RN"
$ws.Range("V8").Value = "117,50"
$ws.Range("W8").Value = "12,00"
$ws.Range("X8").Value = "33,00"
$ws.Range("Y8").Value = "'1"
$ws.Range("Z8").Value = "5.829,00"
$ws.Range("AA8").Value = "https://www.braspress.com.br/w/tracking/search?cnpj=10285063000276&documentType=CONHECIMENTO&numero=3755908"

# Row 9
$ws.Range("A9").Value = "153161, "
$ws.Range("B9").Value = "'3755909"
$ws.Range("C9").Value = " - "
$ws.Range("D9").Value = "'10285063000195"
$ws.Range("E9").Value = "MJ COMERCIO E SERVICOS DE INFORMATICA E"
$ws.Range("F9").Value = "'10483444002556"
$ws.Range("G9").Value = "BETANIA LACTEOS S.A."
$ws.Range("H9").Value = "'05/11/2025"
$ws.Range("I9").Value = "DIGITADO"
$ws.Range("J9").Value = "'11/11/2025"
$ws.Range("K9").Value = "'11/11/2025"
$ws.Range("L9").Value = "AWB EM VIAGEM"
$ws.Range("M9").Value = "LIBERACAO DE OCORR. NA EMISSAO"
$ws.Range("N9").Value = "'06/11/2025"
$ws.Range("O9").Value = "RODO"
$ws.Range("P9").Value = "FOR"
$ws.Range("Q9").Value = "FORTALEZA"
$ws.Range("R9").Value = "CE"
$ws.Range("S9").Value = "SSA"
$ws.Range("T9").Value = "SALVADOR"
$ws.Range("U9").Value = "BA"
$ws.Range("V9").Value = "260,51"
$ws.Range("W9").Value = "12,00"
$ws.Range("X9").Value = "48,72"
$ws.Range("Y9").Value = "'3"
$ws.Range("Z9").Value = "12.576,00"
$ws.Range("AA9").Value = "https://www.braspress.com.br/w/tracking/search?cnpj=10285063000195&documentType=CONHECIMENTO&numero=3755909"

# Row 10
$ws.Range("A10").Value = "152553, "
$ws.Range("B10").Value = "'3740043"
$ws.Range("C10").Value = " - "
$ws.Range("D10").Value = "'10285063000195"
$ws.Range("E10").Value = "MJ COMERCIO E SERVICOS DE INFORMATICA E"
$ws.Range("F10").Value = "'2519126000100"
$ws.Range("G10").Value = "CYBERMAX COMPUTADORES LTDA"
$ws.Range("H10").Value = "27/10/2025"
$ws.Range("I10").Value = "DANFE"
$ws.Range("J10").Value = "'05/11/2025"
$ws.Range("K10").Value = "'05/11/2025"
$ws.Range("L10").Value = "FINALIZADO"
$ws.Range("M10").Value = "ENTREGA REALIZADA"
$ws.Range("N10").Value = "'04/11/2025"
$ws.Range("O10").Value = "RODO"
$ws.Range("P10").Value = "FOR"
$ws.Range("Q10").Value = "FORTALEZA"
$ws.Range("R10").Value = "CE"
$ws.Range("S10").Value = "CCT"
$ws.Range("T10").Value = "SAO PAULO"
$ws.Range("U10").Value = "SP"
$ws.Range("V10").Value = "90,01"
$ws.Range("W10").Value = "12,00"
$ws.Range("X10").Value = "11,00"
$ws.Range("Y10").Value = "'3"
$ws.Range("Z10").Value = "532,52"
$ws.Range("AA10").Value = "https://www.braspress.com.br/w/tracking/search?cnpj=10285063000195&documentType=CONHECIMENTO&numero=3740043"

# --- Total de Conhecimentos line moves from row 10 to row 11 ---
$ws.Range("A11").Value = "Total de Conhecimentos: 3"

# --- Parâmetros / footer block shifts down by one row (13-23 -> 14-24) ---
$ws.Range("A13").Value = ""
$ws.Range("A14").Value = "Parâmetros:"
$ws.Range("A15").Value = "Status: Todos"
$ws.Range("A16").Value = "Pesquisar como Grupo de Empresa: Sim"
$ws.Range("A17").Value = "UF Destino: Todos"
$ws.Range("A18").Value = "Modal: Rodoviário"
$ws.Range("A19").Value = "Pesquiar Por: Remetente"

$ws.Range("A21").Value = ""
$ws.Range("A22").Value = "Observação:"
$ws.Range("A23").Value = "Utiliza dados, referente às emissões efetivadas até o momento da extração deste relatório."

# --- New row 24 (new Data: line; extends the sheet's used range/dimension) ---
$ws.Range("A24").Value = "Data: 09/11/2025 12:49"
